# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" worksheet (same column layout as "2021-Q4")
# between the existing "2021-Q4" and "总计" sheets, populates it with the
# single new holding row, and updates the "总计" (totals) sheet with a new
# summary row for 2022-Q1 (kept above the pre-existing 2021-Q4 summary row).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q1" sheet right after "2021-Q4"
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Add($null, $q4)
$q1.Name = "2022-Q1"

# Pull over the header/row-label formatting (bold + bordered style) from
# the "2021-Q4" sheet so the new sheet matches the existing look exactly.
$q4.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

$q4.Range("A2").Copy()
$q1.Range("A2").PasteSpecial(-4122)

# Header row
$q1.Cells.Item(1,2).Value = "基金代码"
$q1.Cells.Item(1,3).Value = "基金名称"
$q1.Cells.Item(1,4).Value = "基金规模"
$q1.Cells.Item(1,5).Value = "股票总仓位"
$q1.Cells.Item(1,6).Value = "仓位占比"
$q1.Cells.Item(1,7).Value = "持有市值(亿元)"
$q1.Cells.Item(1,8).Value = "仓位排名"

# Data row (index 0)
$q1.Cells.Item(2,1).Value = 0

$q1.Cells.Item(2,2).NumberFormat = "@"
$q1.Cells.Item(2,2).Value = "003655"
$q1.Cells.Item(2,2).Style = "Normal"

$q1.Cells.Item(2,3).Value = "信达澳银新财富灵活配置混合"

$q1.Cells.Item(2,4).NumberFormat = "@"
$q1.Cells.Item(2,4).Value = "11.86"
$q1.Cells.Item(2,4).Style = "Normal"

$q1.Cells.Item(2,5).NumberFormat = "@"
$q1.Cells.Item(2,5).Value = "25.86"
$q1.Cells.Item(2,5).Style = "Normal"

$q1.Cells.Item(2,6).NumberFormat = "@"
$q1.Cells.Item(2,6).Value = "0.64"
$q1.Cells.Item(2,6).Style = "Normal"

$q1.Cells.Item(2,7).NumberFormat = "@"
$q1.Cells.Item(2,7).Value = "0.0759"
$q1.Cells.Item(2,7).Style = "Normal"

$q1.Cells.Item(2,8).Value = 9

# ---------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: add a 2022-Q1 row ahead of the
#    existing 2021-Q4 row.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Make room for the new row by duplicating the row-label style onto A3,
# then fill both rows with their final values (simplest & most robust
# way to land on the target end-state).
$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial(-4122)

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q1"
$total.Cells.Item(2,3).Value = 1
$total.Cells.Item(2,4).Value = 0.08

$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(3,2).Value = "2021-Q4"
$total.Cells.Item(3,3).Value = 4
$total.Cells.Item(3,4).Value = 0.62
